# Adds a new "Page2" worksheet (after TestSheet1) with numeric, date and
# text data, mirroring the "Exc_Col" header row from TestSheet1, and fixes
# up the selection state on both sheets to match what Excel leaves behind
# once Page2 becomes the active/selected sheet.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Preserve the A1:D5 selection on the original sheet before we move away
# from it (Excel keeps a sheet's last selection cached even after it stops
# being the active tab).
[void]$ws1.Range("A1:D5").Select()

# New sheet goes right after TestSheet1.
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Page2"

# Header row - reuse the same column headers as TestSheet1.
$ws2.Range("A1").Value = "Exc_Col1"
$ws2.Range("B1").Value = "Exc_Col2"
$ws2.Range("C1").Value = "Exc_Col3"
$ws2.Range("D1").Value = "Exc_Col4"

# Data rows: col A = plain integer, col B = date serial (formatted as a
# short date), col C/D = text labels A1..A4.
$rows = @(
    @(11, 43497, "A1"),
    @(22, 43498, "A2"),
    @(33, 43499, "A3"),
    @(44, 43500, "A4")
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $ws2.Range("A$r").Value = $rows[$i][0]
    $ws2.Range("B$r").Value = $rows[$i][1]
    $ws2.Range("C$r").Value = $rows[$i][2]
    $ws2.Range("D$r").Value = $rows[$i][2]
}

# Apply a short-date display format to the date column. Format the first
# cell, then copy/paste-special its formatting onto the rest so every
# date cell shares a single style entry instead of one xf per cell.
$ws2.Range("B2").NumberFormat = "mm-dd-yy"
[void]$ws2.Range("B2").Copy()
[void]$ws2.Range("B3:B5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Leave the new sheet active, with C4 selected - matches the saved state
# captured in the target workbook.
[void]$ws2.Range("C4").Select()
